$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Model 1 - MaskRCNN): "Not completed " -> "Work in Progress"
$ws.Range("I9").Value = "Work in Progress"

# Row 11 (Model 2 - Base: ResNet...): "Not completed " -> "Work in Progress "
$ws.Range("I11").Value = "Work in Progress "

# Row 13 (Model 3 - Baseline Model): "Completed but poor scores" -> "Completed"
$ws.Range("I13").Value = "Completed"

# Row 15 (Model 4 - YOLOV3): "Completed but poor scores" -> "Completed "
$ws.Range("I15").Value = "Completed "
# "IoU:.75`nAccuracy ?:Cannot read" -> "IoU:.75`nAccuracy ?"
$ws.Range("H15").Value = "IoU:.75`nAccuracy ?"

# Row 17 (Model 5 - YOLOV4): "Not completed " -> "Work in Progress"
$ws.Range("I17").Value = "Work in Progress"

# Row 19 (Model 6 - SSD): add Accuracy/IOU Score and One Line Status
$ws.Range("H19").Value = "IoU:.75`nAccuracy ?"
$ws.Range("I19").Value = "Completed"
$ws.Range("H19").WrapText = $true
$ws.Rows("19:19").RowHeight = 30
